$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new time-registration row (row 7) that was previously an
# empty placeholder row. Set the number format on F7 to match the other
# "Estimeret tidsforbrug" cells (h:mm) before writing values so the cell
# picks up the existing time-formatted style.
$ws.Range("F7").NumberFormat = "h:mm"
$ws.Range("A7").Value = "Lavet OC0101"
$ws.Range("B7").Value = "System Analyst "
$ws.Range("C7").Value = 43886
$ws.Range("D7").Value = 0.56944444444444442
$ws.Range("E7").Value = 0.59722222222222221
$ws.Range("F7").Value = 0.041666666666666664

# Rename the sheet title (this is the only cell that used the old first
# shared string, so retitling it lets that stale string drop out and the
# workbook's shared-string table renumber itself around the new entries).
$ws.Range("A1").Value = "Tidsregistrering af Emil"

# Move the active selection to D14 (no more frozen/scrolled topLeftCell).
[void]$ws.Range("D14").Select()
